# Adding lab 22 and 23
# Updates the "Asphyxia Test" data table on Sheet1 with the lab-23 values:
# new Blood Pressure / Blood pH readings and refreshed measurement numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Blood Pressure(mmHg): Control / 30 Sec / 1 Min
$ws.Range("B3").Value = "124/79"
$ws.Range("C3").Value = "227/194"
$ws.Range("D3").Value = "64/56"

# Row 4 - Cardiac Output(mL/min)
$ws.Range("B4").Value = 5346
$ws.Range("C4").Value = 10766
$ws.Range("D4").Value = 2574

# Row 6 - Stroke Volume(mL)
$ws.Range("B6").Value = 75
$ws.Range("C6").Value = 55
$ws.Range("D6").Value = 13

# Row 8 - Blood pH(unitless) Arterial/Venous readings
$ws.Range("B8").Value = "7.43/7.38"
$ws.Range("C8").Value = "7.38/7.37"
$ws.Range("D8").Value = "7.37/7.37"
$ws.Range("E8").Value = "7.41/7.41"
